# Weekly data refresh: insert two new daily-price rows (one "Primera" and
# one "Segunda" quality record) for Vega Central Mapocho de Santiago -
# Pepino ensalada, ahead of the existing historical rows, shifting the
# rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 142 (pushes old rows 142:170 down to 144:172).
$ws.Rows.Item(142).Insert()
$ws.Rows.Item(142).Insert()

# --- New row 142: "Primera" quality record ---
$ws.Cells.Item(142, 1).Value = 9
$ws.Cells.Item(142, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(142, 3).Value = "Metropolitana"
$ws.Cells.Item(142, 4).Value = 44508
$ws.Cells.Item(142, 5).Value = 13
$ws.Cells.Item(142, 6).Value = 100112043
$ws.Cells.Item(142, 7).Value = "Pepino ensalada"
$ws.Cells.Item(142, 8).Value = "Sin especificar"
$ws.Cells.Item(142, 9).Value = "Primera"
$ws.Cells.Item(142, 10).Value = 124
$ws.Cells.Item(142, 11).Value = 7000
$ws.Cells.Item(142, 12).Value = 8000
$ws.Cells.Item(142, 13).Value = 7500
$ws.Cells.Item(142, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(142, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(142, 16).Value = 125
$ws.Cells.Item(142, 17).Value = 60
$ws.Cells.Item(142, 18).Value = "Hortaliza"

# --- New row 143: "Segunda" quality record (same date) ---
$ws.Cells.Item(143, 1).Value = 9
$ws.Cells.Item(143, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(143, 3).Value = "Metropolitana"
$ws.Cells.Item(143, 4).Value = 44508
$ws.Cells.Item(143, 5).Value = 13
$ws.Cells.Item(143, 6).Value = 100112043
$ws.Cells.Item(143, 7).Value = "Pepino ensalada"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Segunda"
$ws.Cells.Item(143, 10).Value = 79
$ws.Cells.Item(143, 11).Value = 6000
$ws.Cells.Item(143, 12).Value = 6000
$ws.Cells.Item(143, 13).Value = 6000
$ws.Cells.Item(143, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(143, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(143, 16).Value = 60
$ws.Cells.Item(143, 17).Value = 100
$ws.Cells.Item(143, 18).Value = "Hortaliza"
